{"js": "// \"Complemento Mec\u00e2nica da batalha\"\n// Appends extra detail to the end of the \"Ao come\u00e7ar a partida...\" paragraph\n// (end of the Batalha section) and adds a new closing paragraph describing\n// victory/defeat conditions.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph that talks about starting the match (last paragraph\n// of the \"Batalha\" section in the original document) by matching its\n// distinctive leading text, so the script is resilient to exact indexing.\nconst anchorStart = \"Ao come\u00e7ar a partida\";\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(anchorStart) === 0) {\n    targetParagraph = p;\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not locate the anchor paragraph starting with '\" + anchorStart + \"'.\");\n}\n\n// 1) Append the new sentences about effect/trap cards to the end of that\n//    paragraph (no leading space is added here because the appended text\n//    already starts with one, matching the target paragraph text exactly).\nconst additionalSentences =\n  \" As cartas de efeito servem para causar algum efeito em algum componente do jogo (por exemplo, aumentar o dano que os fighters causam no inimigo). \" +\n  \"J\u00e1 as cartas armadilha ficam esperando que algum jogador inimigo passe por cima dela, causando assim algum efeito no mesmo (como reduzir a velocidade da movimenta\u00e7\u00e3o ou reduzir a quantidade de vida por exemplo).\";\n\ntargetParagraph.insertText(additionalSentences, \"End\");\n\n// 2) Insert a brand-new paragraph right after it, describing the win/lose\n//    condition of a battle.\nconst newParagraphText =\n  \"O Jogador dever\u00e1 utilizar as cartas, criando uma estrat\u00e9gia para destruir a base inimiga sem deixar que o oponente destrua sua base. \" +\n  \"Uma partida normal ter\u00e1 fim quando uma das bases for destru\u00edda. Em algumas batalhas, condi\u00e7\u00f5es especiais encerrar\u00e3o a batalha. \" +\n  \"Se a base destru\u00edda for a base inimiga, o jogador ganha o jogo. Caso contr\u00e1rio, ele perde.\";\n\ntargetParagraph.insertParagraph(newParagraphText, \"After\");\n\nawait context.sync();\n", "ps1": "# \"Complemento Mec\u00e2nica da batalha\"\n# Appends extra detail to the end of the \"Ao come\u00e7ar a partida...\" paragraph\n# (end of the Batalha section) and adds a new closing paragraph describing\n# victory/defeat conditions.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (the last paragraph of the \"Batalha\" section)\n# by searching for its distinctive opening text, so the script does not\n# depend on a hard-coded paragraph index.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"Ao come\u00e7ar a partida\")\nif (-not $found) {\n    throw \"Could not locate the anchor paragraph starting with 'Ao come\u00e7ar a partida'.\"\n}\n\n$targetParagraph = $searchRange.Paragraphs(1)\n$targetRange = $targetParagraph.Range\n\n# 1) Append the new sentences about effect/trap cards to the end of that\n#    paragraph (leading space included so it reads naturally after the\n#    existing final sentence).\n$additionalSentences = \" As cartas de efeito servem para causar algum efeito em algum componente do jogo (por exemplo, aumentar o dano que os fighters causam no inimigo). J\u00e1 as cartas armadilha ficam esperando que algum jogador inimigo passe por cima dela, causando assim algum efeito no mesmo (como reduzir a velocidade da movimenta\u00e7\u00e3o ou reduzir a quantidade de vida por exemplo).\"\n$targetRange.InsertAfter($additionalSentences)\n\n# 2) Insert a brand-new paragraph right after it, describing the win/lose\n#    condition of a battle.\n$targetRange.InsertParagraphAfter()\n$newParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$newParagraphText = \"O Jogador dever\u00e1 utilizar as cartas, criando uma estrat\u00e9gia para destruir a base inimiga sem deixar que o oponente destrua sua base. Uma partida normal ter\u00e1 fim quando uma das bases for destru\u00edda. Em algumas batalhas, condi\u00e7\u00f5es especiais encerrar\u00e3o a batalha. Se a base destru\u00edda for a base inimiga, o jogador ganha o jogo. Caso contr\u00e1rio, ele perde.\"\n$newParagraph.Range.InsertAfter($newParagraphText)\n"}
